# Update "Des Scheduled Flights vs actual.xlsx"
# Extend the flight-tracking table (Ark1) from row 392 through row 406
# with 14 more days of data (2021-05-03 .. 2021-05-16), carrying the
# existing formatting / formula pattern down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting (number formats, fonts, alignment, etc.) of the
#    last existing data row (392) down across the new rows so the new
#    cells pick up the same styles already used by the table (s=1/2/3)
#    instead of Excel inventing brand-new style records.
$ws.Range("A392:D392").Copy()
$ws.Range("A393:D406").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 2) New row data: Date, Scheduled flights (B), Tracked flights (C)
$rows = @(
    @{ Row = 393; Date = "2021-05-03"; B = 50; C = 50 },
    @{ Row = 394; Date = "2021-05-04"; B = 57; C = 52 },
    @{ Row = 395; Date = "2021-05-05"; B = 62; C = 59 },
    @{ Row = 396; Date = "2021-05-06"; B = 64; C = 63 },
    @{ Row = 397; Date = "2021-05-07"; B = 66; C = 60 },
    @{ Row = 398; Date = "2021-05-08"; B = 41; C = 40 },
    @{ Row = 399; Date = "2021-05-09"; B = 48; C = 47 },
    @{ Row = 400; Date = "2021-05-10"; B = 67; C = 62 },
    @{ Row = 401; Date = "2021-05-11"; B = 58; C = 52 },
    @{ Row = 402; Date = "2021-05-12"; B = 65; C = 60 },
    @{ Row = 403; Date = "2021-05-13"; B = 66; C = 66 },
    @{ Row = 404; Date = "2021-05-14"; B = 65; C = 61 },
    @{ Row = 405; Date = "2021-05-15"; B = 50; C = 47 },
    @{ Row = 406; Date = "2021-05-16"; B = 53; C = 52 }
)

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value2 = $r.Date
    $ws.Range("B$($r.Row)").Value2 = $r.B
    $ws.Range("C$($r.Row)").Value2 = $r.C
}

# 3) Fill column D with the Tracked/Scheduled ratio formula, extending the
#    existing pattern (=C{row}/B{row}) down through row 406 in one shot so
#    the engine groups it as a shared formula, matching how the original
#    author dragged the formula down.
$ws.Range("D393:D406").Formula = "=C393/B393"

# 4) Update the view: the sheet had scrolled further down and the active
#    selection moved to the newly appended D column range.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 323
$aw.ScrollColumn = 1
$ws.Range("D392:D406").Select()
